$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country/timestamp labels (column A) ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Marzo de 2020 a las 15:46"
$ws.Cells.Item(72, 1).Value = "Republica Dominicana"
$ws.Cells.Item(73, 1).Value = "Uruguay"
$ws.Cells.Item(74, 1).Value = "Hungria"
$ws.Cells.Item(75, 1).Value = "Argelia"
$ws.Cells.Item(76, 1).Value = "Vietnam"
$ws.Cells.Item(77, 1).Value = "Islas Feroe"
$ws.Cells.Item(78, 1).Value = "Bosnia y Herzegovina"
$ws.Cells.Item(79, 1).Value = "Principado de Andorra"
$ws.Cells.Item(80, 1).Value = "Marruecos"
$ws.Cells.Item(82, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(83, 1).Value = "Brunei"
$ws.Cells.Item(84, 1).Value = "Sri Lanka"
$ws.Cells.Item(85, 1).Value = "Albania"
$ws.Cells.Item(86, 1).Value = "Bielorrusia"
$ws.Cells.Item(87, 1).Value = "Republica de Chipre"
$ws.Cells.Item(88, 1).Value = "Malta"
$ws.Cells.Item(95, 1).Value = "Camboya"
$ws.Cells.Item(96, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(97, 1).Value = "Oman"
$ws.Cells.Item(98, 1).Value = "Estado de Palestina"
$ws.Cells.Item(115, 1).Value = "Bolivia"
$ws.Cells.Item(116, 1).Value = "Ghana"
$ws.Cells.Item(119, 1).Value = "Puerto Rico"
$ws.Cells.Item(120, 1).Value = "Ruanda"
$ws.Cells.Item(121, 1).Value = "Macao"
$ws.Cells.Item(122, 1).Value = "Togo"
$ws.Cells.Item(123, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(124, 1).Value = "Guam"
$ws.Cells.Item(125, 1).Value = "Guayana Francesa"
$ws.Cells.Item(126, 1).Value = "Kirguistan"
$ws.Cells.Item(127, 1).Value = "Montenegro"
$ws.Cells.Item(135, 1).Value = "Etiopia"
$ws.Cells.Item(136, 1).Value = "Trinidad yTobago"
$ws.Cells.Item(137, 1).Value = "Mayotte"
$ws.Cells.Item(138, 1).Value = "Seychelles"
$ws.Cells.Item(139, 1).Value = "Kenia"
$ws.Cells.Item(140, 1).Value = "Barbados"
$ws.Cells.Item(141, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(142, 1).Value = "Tanzania"
$ws.Cells.Item(143, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(146, 1).Value = "San Martin (Parte Francesa)"
$ws.Cells.Item(148, 1).Value = "Surinam"
$ws.Cells.Item(151, 1).Value = "Cabo Verde"
$ws.Cells.Item(152, 1).Value = "Congo"
$ws.Cells.Item(154, 1).Value = "Namibia"
$ws.Cells.Item(155, 1).Value = "San Bartolome"
$ws.Cells.Item(156, 1).Value = "Zimbabue"
$ws.Cells.Item(157, 1).Value = "Republica de Africa Central"
$ws.Cells.Item(159, 1).Value = "Madagascar"
$ws.Cells.Item(162, 1).Value = "Mauritania"
$ws.Cells.Item(163, 1).Value = "Isla de Man"
$ws.Cells.Item(164, 1).Value = "Haiti"
$ws.Cells.Item(165, 1).Value = "Guinea"
$ws.Cells.Item(166, 1).Value = "Benin"
$ws.Cells.Item(167, 1).Value = "Santa Lucia"
$ws.Cells.Item(168, 1).Value = "Angola"
$ws.Cells.Item(169, 1).Value = "Butan"
$ws.Cells.Item(170, 1).Value = "Zambia"
$ws.Cells.Item(171, 1).Value = "Nicaragua"
$ws.Cells.Item(173, 1).Value = "Bermudas"
$ws.Cells.Item(174, 1).Value = "Fiyi"
$ws.Cells.Item(176, 1).Value = "Suazilandia"
$ws.Cells.Item(177, 1).Value = "Montserrat"
$ws.Cells.Item(178, 1).Value = "Niger"
$ws.Cells.Item(179, 1).Value = "Republica del Chad"
$ws.Cells.Item(180, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(181, 1).Value = "Timor Oriental"
$ws.Cells.Item(183, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(184, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(185, 1).Value = "Somalia"
$ws.Cells.Item(186, 1).Value = "Gambia"
$ws.Cells.Item(187, 1).Value = "San Martin (Parte Holandesa)"
$ws.Cells.Item(188, 1).Value = "Antigua y Barbuda"

# --- Update statistic values (columns B-H) ---
$ws.Cells.Item(7, 2).Value = 21652
$ws.Cells.Item(7, 3).Value = 1804
$ws.Cells.Item(7, 5).Value = 21370
$ws.Cells.Item(9, 2).Value = 19823
$ws.Cells.Item(9, 3).Value = 440
$ws.Cells.Item(9, 5).Value = 19400
$ws.Cells.Item(13, 2).Value = 4072
$ws.Cells.Item(13, 3).Value = 89
$ws.Cells.Item(13, 5).Value = 3830
$ws.Cells.Item(17, 2).Value = 2017
$ws.Cells.Item(17, 3).Value = 58
$ws.Cells.Item(17, 5).Value = 2009
$ws.Cells.Item(25, 2).Value = 998
$ws.Cells.Item(25, 3).Value = 28
$ws.Cells.Item(25, 5).Value = 984
$ws.Cells.Item(33, 4).Value = 8
$ws.Cells.Item(33, 5).Value = 529
$ws.Cells.Item(47, 4).Value = 2
$ws.Cells.Item(47, 5).Value = 304
$ws.Cells.Item(63, 5).Value = 151
$ws.Cells.Item(63, 7).Value = 1
$ws.Cells.Item(63, 8).Value = 4
$ws.Cells.Item(65, 4).Value = 38
$ws.Cells.Item(65, 5).Value = 113
$ws.Cells.Item(72, 2).Value = 112
$ws.Cells.Item(72, 3).Value = 40
$ws.Cells.Item(72, 8).Value = 2
$ws.Cells.Item(73, 2).Value = 110
$ws.Cells.Item(73, 3).Value = 0
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 110
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 8).Value = 0
$ws.Cells.Item(74, 2).Value = 103
$ws.Cells.Item(74, 3).Value = 18
$ws.Cells.Item(74, 4).Value = 7
$ws.Cells.Item(74, 5).Value = 92
$ws.Cells.Item(74, 6).Value = 6
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 4
$ws.Cells.Item(75, 2).Value = 95
$ws.Cells.Item(75, 3).Value = 1
$ws.Cells.Item(75, 4).Value = 43
$ws.Cells.Item(75, 5).Value = 40
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 1
$ws.Cells.Item(75, 8).Value = 12
$ws.Cells.Item(76, 2).Value = 94
$ws.Cells.Item(76, 3).Value = 3
$ws.Cells.Item(76, 4).Value = 17
$ws.Cells.Item(76, 5).Value = 77
$ws.Cells.Item(76, 6).Value = 2
$ws.Cells.Item(77, 2).Value = 92
$ws.Cells.Item(77, 3).Value = 12
$ws.Cells.Item(77, 4).Value = 3
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(78, 2).Value = 91
$ws.Cells.Item(78, 3).Value = 2
$ws.Cells.Item(78, 4).Value = 2
$ws.Cells.Item(78, 5).Value = 89
$ws.Cells.Item(78, 6).Value = 1
$ws.Cells.Item(79, 2).Value = 88
$ws.Cells.Item(79, 3).Value = 13
$ws.Cells.Item(79, 4).Value = 1
$ws.Cells.Item(79, 5).Value = 87
$ws.Cells.Item(79, 6).Value = 2
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(80, 2).Value = 86
$ws.Cells.Item(80, 3).Value = 0
$ws.Cells.Item(80, 4).Value = 2
$ws.Cells.Item(80, 5).Value = 81
$ws.Cells.Item(80, 8).Value = 3
$ws.Cells.Item(82, 2).Value = 85
$ws.Cells.Item(82, 3).Value = 9
$ws.Cells.Item(82, 5).Value = 84
$ws.Cells.Item(82, 6).Value = 1
$ws.Cells.Item(83, 2).Value = 83
$ws.Cells.Item(83, 3).Value = 5
$ws.Cells.Item(83, 4).Value = 1
$ws.Cells.Item(83, 5).Value = 82
$ws.Cells.Item(84, 2).Value = 77
$ws.Cells.Item(84, 3).Value = 4
$ws.Cells.Item(84, 4).Value = 3
$ws.Cells.Item(84, 5).Value = 74
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(85, 3).Value = 6
$ws.Cells.Item(85, 4).Value = 2
$ws.Cells.Item(85, 5).Value = 72
$ws.Cells.Item(85, 6).Value = 2
$ws.Cells.Item(85, 8).Value = 2
$ws.Cells.Item(86, 2).Value = 76
$ws.Cells.Item(86, 3).Value = 7
$ws.Cells.Item(86, 4).Value = 15
$ws.Cells.Item(86, 5).Value = 61
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(87, 2).Value = 75
$ws.Cells.Item(87, 3).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 75
$ws.Cells.Item(88, 2).Value = 73
$ws.Cells.Item(88, 3).Value = 9
$ws.Cells.Item(88, 4).Value = 2
$ws.Cells.Item(88, 5).Value = 71
$ws.Cells.Item(88, 6).Value = 1
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(95, 2).Value = 53
$ws.Cells.Item(95, 3).Value = 2
$ws.Cells.Item(95, 4).Value = 1
$ws.Cells.Item(96, 3).Value = 0
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 52
$ws.Cells.Item(97, 4).Value = 13
$ws.Cells.Item(97, 5).Value = 39
$ws.Cells.Item(98, 2).Value = 52
$ws.Cells.Item(98, 3).Value = 4
$ws.Cells.Item(98, 4).Value = 17
$ws.Cells.Item(98, 5).Value = 35
$ws.Cells.Item(100, 2).Value = 49
$ws.Cells.Item(100, 3).Value = 5
$ws.Cells.Item(100, 5).Value = 48
$ws.Cells.Item(119, 3).Value = 3
$ws.Cells.Item(120, 4).Value = 0
$ws.Cells.Item(120, 5).Value = 17
$ws.Cells.Item(121, 2).Value = 17
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(121, 4).Value = 10
$ws.Cells.Item(121, 5).Value = 7
$ws.Cells.Item(122, 2).Value = 16
$ws.Cells.Item(122, 3).Value = 7
$ws.Cells.Item(122, 5).Value = 16
$ws.Cells.Item(123, 3).Value = 4
$ws.Cells.Item(124, 3).Value = 1
$ws.Cells.Item(125, 2).Value = 15
$ws.Cells.Item(125, 5).Value = 15
$ws.Cells.Item(126, 3).Value = 8
$ws.Cells.Item(127, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 3
$ws.Cells.Item(131, 5).Value = 10
$ws.Cells.Item(141, 3).Value = 0
$ws.Cells.Item(143, 3).Value = 3
$ws.Cells.Item(151, 3).Value = 2
$ws.Cells.Item(152, 3).Value = 0
$ws.Cells.Item(156, 3).Value = 2
$ws.Cells.Item(159, 3).Value = 0
$ws.Cells.Item(162, 3).Value = 0
$ws.Cells.Item(168, 3).Value = 1
$ws.Cells.Item(173, 3).Value = 0
$ws.Cells.Item(174, 3).Value = 1
$ws.Cells.Item(181, 3).Value = 1
$ws.Cells.Item(185, 3).Value = 0
